# edit.ps1
# Applies the diff: adds en-US language markers to the first three paragraphs,
# inserts a new block of paragraphs (form fields 5,6,8 with checkboxes) after
# the "Antrag" paragraph, and adds a new "hljs-number" character style.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Add <w:lang w:val="en-US"/> to paragraph mark (pPr/rPr) and run (r/rPr)
#        of the first three paragraphs ({#rows}, Account Description, Antrag) ---

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{#rows}</w:t></w:r></w:p>')

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Account Description: {AccountDescription}</w:t></w:r></w:p>')

$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Antrag: {Antrag}</w:t></w:r></w:p>')

# --- 2. Insert the new block of paragraphs (category / competence / delegation
#        checkboxes) right after the "Antrag" paragraph, before the dashed line ---

$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>5)  Kategorie gem. &#167; 42 Abs. 2 BO:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">{ID_1} </w:t>
  </w:r>
  <w:r>
    <w:t>Anforderung, die bereits in der Investitionsvorschau (&#167; 44) ber&#252;cksichtigt ist</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">{ID_2} </w:t>
  </w:r>
  <w:r>
    <w:t>Anforderung au&#223;erhalb der Investitionsvorschau, ausgenommen laufender Bedarf</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">{ID_3} </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Laufender </w:t>
  </w:r>
  <w:r>
    <w:t>Bedarf</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t xml:space="preserve">6)  Sachliche Zust&#228;ndigkeit zur Bewilligung gem. &#167; 43 Abs. 1 bzw. Abs. 3 BO): </w:t>
  </w:r>
  <w:r>
    <w:br/>
  </w:r>
  <w:r>
    <w:t>{ID_</w:t>
  </w:r>
  <w:r>
    <w:t>4</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">} </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Bau, Facility Management und Beschaffung</w:t>
  </w:r>
  <w:r>
    <w:br/>
  </w:r>
  <w:r>
    <w:t>{ID_</w:t>
  </w:r>
  <w:r>
    <w:t>5</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">} </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">  IT Management und IT Organisation</w:t>
  </w:r>
  <w:r>
    <w:br/>
  </w:r>
  <w:r>
    <w:t>{ID_</w:t>
  </w:r>
  <w:r>
    <w:t>6</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">} </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> &#214;A, Marketing</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:spacing w:after="0"/>
  </w:pPr>
  <w:r>
    <w:t>8)  Pr&#252;fung der Delegierung an das B&#252;ro:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:left="1440" w:hanging="720"/>
  </w:pPr>
  <w:r>
    <w:t>{ID_</w:t>
  </w:r>
  <w:r>
    <w:t>7</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">} </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> B&#252;ro</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>gem. &#167; 41 Abs. 1 Z1 BO in folgenden F&#228;llen</w:t>
  </w:r>
  <w:r>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> &lt; &#8364; 214.500,00  laufende Verwaltungsgesch&#228;fte, soweit im Einzelfall das Eineinhalbfache des f&#252;r das jeweilige Jahr festgesetzten Schwellenwertes f&#252;r Dienstleistungen nach &#167; 12 Abs. 1 Z 1 BVerG 2018 nicht &#252;berschritten wird (C1 AnhGOV)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:left="1440" w:hanging="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:r>
    <w:t>{ID_</w:t>
  </w:r>
  <w:r>
    <w:t>8</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">} </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> bis &#8364; 54.540,00 (Leiterin bzw. Leiter d. &#246;rtl. zust. OE)</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> &#8364; 54.540,01 bis &#8364; 109.080,00 (FB/EZ-Leiterin bzw. FB/EZ-Leiter)</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> &#8364; 109.080,01 bis &#8364; 214.500,00 (Generaldirektor bzw. GB-Leiter)</w:t>
  </w:r>
  <w:r>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Vergabe notwendiger wiederkehrender Auftr&#228;ge zur Sicherstellung eines reibungslosen Betriebsablaufes und zur Aufrechterhaltung der Betriebssicherheit (C13 AnhGOV)</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> bis &#8364; 54.540,00 (Leiterin bzw. Leiter d. &#246;rtl. zust. OE)</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> &#8364; 54.540,01 bis &#8364; 109.080,00 (FB/EZ-Leiterin bzw. FB/EZ-Leiter)</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> ab &#8364; 109.080,01 (Generaldirektor bzw. GB-Leiter)</w:t>
  </w:r>
  <w:r>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Nachschaffung von Gegenst&#228;nden und Materialien sowie f&#252;r notwendige &#220;berpr&#252;fungen und Wartungen, die zur Aufrechterhaltung einer kontinuierlichen Betriebsf&#252;hrung unbedingt erforderlich sind, ohne betragliche Obergrenze (AnhGOVR Abschnitt C, Punkt 12 und 13).</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Wingdings" w:eastAsia="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings"/>
    </w:rPr>
    <w:t>&#xF0A8;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> unbegrenzt (Leiterin bzw. Leiter d. &#246;rtl. zust. OE)</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p/>
'@

$p3 = $d.Paragraphs.Item(3)
$insertionPoint = $d.Range($p3.Range.End, $p3.Range.End)
$insertionPoint.InsertXML($newBlockXml)

# --- 3. Add the new "hljs-number" character style (based on Default Paragraph Font) ---

$newStyle = $d.Styles.Add("hljs-number", 2)
$newStyle.BaseStyle = "DefaultParagraphFont"
